$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.395.15'
$ws.Range('E2').Value = '  +2.67%  '

$ws.Range('D3').Value = '2.695.52'
$ws.Range('E3').Value = '  +1.53%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').Value = '''524.35'
$ws.Range('E5').Value = '  +1.97%  '

$ws.Range('D6').Value = '''147.25'
$ws.Range('E6').Value = '  +2.40%  '

$ws.Range('E7').Value = '  +0.11%  '

$ws.Range('E8').Value = '  +1.96%  '

$ws.Range('D9').Value = '2.717.84'
$ws.Range('E9').Value = '  +1.32%  '

$ws.Range('D10').Value = '''6.51'
$ws.Range('E10').Value = '  +4.97%  '

$ws.Range('D11').Value = '''0.107'
$ws.Range('E11').Value = '  +0.39%  '

$ws.Range('E12').Value = '  +1.65%  '

$ws.Range('E13').Value = '  +1.62%  '

$ws.Range('D14').Value = '3.170.38'
$ws.Range('E14').Value = '  +1.55%  '

$ws.Range('D15').Value = '60.400.79'
$ws.Range('E15').Value = '  +2.72%  '

$ws.Range('D16').Value = '2.847.24'
$ws.Range('E16').Value = '  +6.22%  '

$ws.Range('D17').Value = '''21.45'
$ws.Range('E17').Value = '  +2.34%  '

$ws.Range('E18').Value = '  +1.45%  '

$ws.Range('D19').Value = '''351.96'
$ws.Range('E19').Value = '  +1.58%  '

$ws.Range('E20').Value = '  +0.57%  '

$ws.Range('E21').Value = '  +2.09%  '

$ws.Range('D22').Value = '''6.37'
$ws.Range('E22').Value = '  +3.78%  '

$ws.Range('E23').Value = '  +0.09%  '

$ws.Range('D24').Value = '''63.15'
$ws.Range('E24').Value = '  +3.70%  '

$ws.Range('E25').Value = '  +1.24%  '

$ws.Range('E26').Value = '  +5.68%  '

$ws.Range('D27').Value = '''0.993'
$ws.Range('E27').Value = '  +0.00%  '

$ws.Range('D28').Value = '0.0₃0823'
$ws.Range('E28').Value = '  +1.86%  '

$ws.Range('D29').Value = '''7.37'
$ws.Range('E29').Value = '  +2.18%  '

$ws.Range('D30').Value = '''6.89'
$ws.Range('E30').Value = '  +7.13%  '

$ws.Range('E31').Value = '  +0.10%  '

$ws.Range('E32').Value = '  +1.90%  '

$ws.Range('D33').Value = '''19.19'
$ws.Range('E33').Value = '  +1.28%  '

$ws.Range('E34').Value = '  -1.47%  '

$ws.Range('D35').Value = '''4.31'
$ws.Range('E35').Value = '  +7.13%  '

$ws.Range('E36').Value = '  +10.15%  '

$ws.Range('D37').Value = '''0.961'
$ws.Range('E37').Value = '  -4.75%  '

$ws.Range('E38').Value = '  +9.35%  '

$ws.Range('E39').Value = '  +4.52%  '

$ws.Range('E41').Value = '  +0.84%  '

$ws.Range('D42').Value = '''288.12'
$ws.Range('E42').Value = '  +3.43%  '

$ws.Range('D43').Value = '''20.17'
$ws.Range('E43').Value = '  +2.21%  '

$ws.Range('E44').Value = '  -0.29%  '

$ws.Range('D45').Value = '''0.0993'
$ws.Range('E45').Value = '  +1.34%  '

$ws.Range('E46').Value = '  +0.35%  '

$ws.Range('D47').Value = '2.144.03'
$ws.Range('E47').Value = '  +7.02%  '

$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').Value = '''0.0542'
$ws.Range('E48').Value = '  +1.77%  '

$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '''4.92'
$ws.Range('E49').Value = '  +4.08%  '

$ws.Range('E50').Value = '  +2.35%  '

$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '''19.37'
$ws.Range('E51').Value = '  +6.73%  '
